$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32 (shifts existing rows 32-44 down to 33-45)
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the latest weekly record
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44636
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112031
$ws.Range("G32").Value = "Poroto verde"
$ws.Range("H32").Value = "Magnum"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 180
$ws.Range("K32").Value = 22000
$ws.Range("L32").Value = 23000
$ws.Range("M32").Value = 22444
$ws.Range("N32").Value = "$/saco 25 kilos"
$ws.Range("O32").Value = "Región del Maule"
$ws.Range("P32").Value = 898
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
